$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    ,(0.081928, 0.02525817399873631)
    ,(0.163856, 0.02543574199808063)
    ,(0.245784, 0.02524120299494825)
    ,(0.901209, 0.02510983100364683)
    ,(0.983137, 0.02586906700162217)
    ,(1.065065, 0.02542278699547751)
    ,(1.310849, 0.02564401899871882)
    ,(1.392777, 0.02543641599913826)
    ,(1.884346, 0.02434378200268839)
    ,(2.13013, 0.0250395330003812)
    ,(2.212058, 0.02640118500130484)
    ,(2.375914, 0.02666043900535442)
    ,(2.457842, 0.02645270300126867)
    ,(2.621698, 0.02654777800489683)
    ,(2.703627, 0.02665507500205422)
    ,(3.031339, 0.02583294400392333)
    ,(3.359051, 0.02537747300084447)
    ,(3.604836, 0.0258221969997976)
    ,(3.768692, 0.02522520600177813)
    ,(3.85062, 0.02532092900219141)
    ,(4.014476, 0.02668907799670706)
    ,(4.506045, 0.02538750600069761)
    ,(4.669901, 0.02610378600365948)
    ,(4.751829, 0.02618399700440932)
    ,(4.997613, 0.02570776799984742)
    ,(5.243397, 0.02638006699999096)
    ,(5.325326, 0.02556784199987305)
    ,(5.57111, 0.02583414399850881)
    ,(5.816894, 0.02595475999987684)
    ,(5.898822, 0.02630334900459275)
    ,(6.062678, 0.0262667130009504)
    ,(6.144607, 0.0259565059968736)
    ,(6.390391, 0.02645851299894275)
    ,(6.472319, 0.02667079800448846)
    ,(6.718103, 0.02650103299674811)
    ,(6.963887, 0.02694696399703389)
    ,(7.127744, 0.02769829500175547)
    ,(7.209672, 0.0267537270046887)
    ,(7.2916, 0.02733118700416526)
    ,(7.537384, 0.02851465599815128)
    ,(7.619312, 0.02950993299600668)
    ,(7.783168, 0.03180854800302768)
    ,(7.865096, 0.03074928199930582)
    ,(8.028953, 0.03147868000087328)
    ,(8.356665, 0.03118221999466186)
    ,(8.438592999999999, 0.02955734500574181)
    ,(8.602449, 0.03141222600243054)
    ,(8.684377, 0.0302616230037529)
    ,(9.503658, 0.03110415799892507)
    ,(9.995227, 0.03078853199986042)
    ,(10.241011, 0.0316529060000903)
    ,(10.322939, 0.03194444199471036)
    ,(10.650652, 0.03093208099744515)
    ,(10.73258, 0.02972358599799918)
    ,(10.896436, 0.03225036799994996)
    ,(12.043429, 0.02924031800648663)
    ,(12.534998, 0.03001804900122806)
    ,(12.616926, 0.03015428299841005)
    ,(12.86271, 0.0281556210029521)
    ,(13.272351, 0.02966268199816113)
    ,(13.436207, 0.02974016799998935)
    ,(13.518135, 0.03018881700700149)
    ,(14.009704, 0.0299486580042867)
    ,(14.091632, 0.03040733800298767)
    ,(14.255488, 0.03009416799613973)
    ,(14.337416, 0.03006054500292521)
    ,(14.419344, 0.03085615799500374)
    ,(14.910913, 0.03079165399685735)
    ,(15.238625, 0.03033495199633762)
    ,(15.730193, 0.03079555999283912)
    ,(15.812122, 0.02989649200026179)
    ,(16.30369, 0.03123082999809412)
    ,(16.385618, 0.03011389399762265)
    ,(16.631403, 0.03096671800449258)
    ,(16.959115, 0.03074913100135745)
    ,(17.122971, 0.03169431500282371)
    ,(17.204899, 0.03118541800358798)
    ,(17.368755, 0.03179089300101623)
    ,(17.450683, 0.03304675900290022)
    ,(17.778396, 0.03275902600580594)
    ,(18.02418, 0.03175470499991206)
    ,(18.188036, 0.0319614459949662)
    ,(18.515749, 0.03142127099999925)
    ,(18.597677, 0.03270616899681045)
    ,(18.761533, 0.03273856099986006)
    ,(18.843461, 0.03356986599828815)
    ,(18.925389, 0.03294246699806536)
    ,(19.416958, 0.03448192699579522)
    ,(19.498886, 0.03377556299528806)
    ,(19.662742, 0.03487447099905694)
    ,(19.74467, 0.03408476799813798)
    ,(19.990454, 0.03873087400279474)
    ,(20.072382, 0.03704707400174811)
    ,(20.482023, 0.04039318600553088)
    ,(20.563951, 0.04061458000069251)
    ,(20.645879, 0.0382942619980895)
    ,(20.809735, 0.04105730000446783)
    ,(20.891663, 0.04221983200113755)
    ,(21.629016, 0.04121946400118759)
    ,(21.710944, 0.04135833100008313)
    ,(21.874801, 0.04122477499913657)
    ,(22.448297, 0.0412707750001573)
    ,(22.530225, 0.04085884999949485)
    ,(22.612153, 0.04125330100214342)
    ,(23.103722, 0.0421513299952494)
    ,(23.18565, 0.04345248600293417)
    ,(23.677219, 0.04343833099846961)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
